# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el resumen de conversión del día ---
$wsResumen = $wb.Worksheets.Item("Hoja1")

$texto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.64 = 26317.36 pesos`n✅ 26317.36 pesos = 6.65 = 975.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsResumen.Range("A1").Value = $texto

# --- tasas: actualizar las tasas de Binance y transfi ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 150.509
$wsTasas.Range("O10").Value = 3961
$wsTasas.Range("N12").Value = 3958
$wsTasas.Range("O12").Value = 146.65
